$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (8 and 9) before populating so existing rows 2-7 keep their positions
$ws.Rows("8:9").Insert()

# Row 2: ECs -> Efnb1/Ephb1 -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Efnb1"
$ws.Cells.Item(2,3).Value = "Ephb1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 11.43712066666667
$ws.Cells.Item(2,8).Value = 34.311362
$ws.Cells.Item(2,9).Value = 0.5796330080444665
$ws.Cells.Item(2,10).Value = 0.5796330080444665
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 1.501929333333333
$ws.Cells.Item(2,14).Value = 4.505788
$ws.Cells.Item(2,15).Value = 0.7650463650777426
$ws.Cells.Item(2,16).Value = 0.7650463650777426
$ws.Cells.Item(2,17).Value = 17.17774701813956
$ws.Cells.Item(2,18).Value = 154.599723163256
$ws.Cells.Item(2,19).Value = 0.443446125883497
$ws.Cells.Item(2,20).Value = 0.443446125883497

# Row 3: ECs -> Efnb1/Ephb1 -> sCs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Efnb1"
$ws.Cells.Item(3,3).Value = "Ephb1"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 11.43712066666667
$ws.Cells.Item(3,8).Value = 34.311362
$ws.Cells.Item(3,9).Value = 0.5796330080444665
$ws.Cells.Item(3,10).Value = 0.5796330080444665
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.461258
$ws.Cells.Item(3,14).Value = 1.383774
$ws.Cells.Item(3,15).Value = 0.2349536349222574
$ws.Cells.Item(3,16).Value = 0.2349536349222574
$ws.Cells.Item(3,17).Value = 5.275463404465333
$ws.Cells.Item(3,18).Value = 47.47917064018801
$ws.Cells.Item(3,19).Value = 0.1361868821609695
$ws.Cells.Item(3,20).Value = 0.1361868821609695

# Row 4: FAPs -> Efnb1/Ephb1 -> ECs
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Efnb1"
$ws.Cells.Item(4,3).Value = "Ephb1"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 5.754308333333334
$ws.Cells.Item(4,8).Value = 17.262925
$ws.Cells.Item(4,9).Value = 0.2916282118266253
$ws.Cells.Item(4,10).Value = 0.2916282118266253
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.501929333333333
$ws.Cells.Item(4,14).Value = 4.505788
$ws.Cells.Item(4,15).Value = 0.7650463650777426
$ws.Cells.Item(4,16).Value = 0.7650463650777426
$ws.Cells.Item(4,17).Value = 8.642564478877778
$ws.Cells.Item(4,18).Value = 77.78308030990001
$ws.Cells.Item(4,19).Value = 0.2231091034120816
$ws.Cells.Item(4,20).Value = 0.2231091034120816

# Row 5: FAPs -> Efnb1/Ephb1 -> sCs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Efnb1"
$ws.Cells.Item(5,3).Value = "Ephb1"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 5.754308333333334
$ws.Cells.Item(5,8).Value = 17.262925
$ws.Cells.Item(5,9).Value = 0.2916282118266253
$ws.Cells.Item(5,10).Value = 0.2916282118266253
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.461258
$ws.Cells.Item(5,14).Value = 1.383774
$ws.Cells.Item(5,15).Value = 0.2349536349222574
$ws.Cells.Item(5,16).Value = 0.2349536349222574
$ws.Cells.Item(5,17).Value = 2.654220753216667
$ws.Cells.Item(5,18).Value = 23.88798677895
$ws.Cells.Item(5,19).Value = 0.06851910841454367
$ws.Cells.Item(5,20).Value = 0.06851910841454367

# Row 6: M2 -> Efnb1/Ephb1 -> ECs
$ws.Cells.Item(6,1).Value = "M2"
$ws.Cells.Item(6,2).Value = "Efnb1"
$ws.Cells.Item(6,3).Value = "Ephb1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.2670106666666667
$ws.Cells.Item(6,8).Value = 0.8010320000000001
$ws.Cells.Item(6,9).Value = 0.0135320943453039
$ws.Cells.Item(6,10).Value = 0.0135320943453039
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.501929333333333
$ws.Cells.Item(6,14).Value = 4.505788
$ws.Cells.Item(6,15).Value = 0.7650463650777426
$ws.Cells.Item(6,16).Value = 0.7650463650777426
$ws.Cells.Item(6,17).Value = 0.4010311525795555
$ws.Cells.Item(6,18).Value = 3.609280373216
$ws.Cells.Item(6,19).Value = 0.01035267959076382
$ws.Cells.Item(6,20).Value = 0.01035267959076382

# Row 7: M2 -> Efnb1/Ephb1 -> sCs
$ws.Cells.Item(7,1).Value = "M2"
$ws.Cells.Item(7,2).Value = "Efnb1"
$ws.Cells.Item(7,3).Value = "Ephb1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.2670106666666667
$ws.Cells.Item(7,8).Value = 0.8010320000000001
$ws.Cells.Item(7,9).Value = 0.0135320943453039
$ws.Cells.Item(7,10).Value = 0.0135320943453039
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.461258
$ws.Cells.Item(7,14).Value = 1.383774
$ws.Cells.Item(7,15).Value = 0.2349536349222574
$ws.Cells.Item(7,16).Value = 0.2349536349222574
$ws.Cells.Item(7,17).Value = 0.1231608060853333
$ws.Cells.Item(7,18).Value = 1.108447254768
$ws.Cells.Item(7,19).Value = 0.003179414754540076
$ws.Cells.Item(7,20).Value = 0.003179414754540076

# Row 8: sCs -> Efnb1/Ephb1 -> ECs
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Efnb1"
$ws.Cells.Item(8,3).Value = "Ephb1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.273219
$ws.Cells.Item(8,8).Value = 6.819656999999999
$ws.Cells.Item(8,9).Value = 0.1152066857836043
$ws.Cells.Item(8,10).Value = 0.1152066857836043
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 1.501929333333333
$ws.Cells.Item(8,14).Value = 4.505788
$ws.Cells.Item(8,15).Value = 0.7650463650777426
$ws.Cells.Item(8,16).Value = 0.7650463650777426
$ws.Cells.Item(8,17).Value = 3.414214297190666
$ws.Cells.Item(8,18).Value = 30.727928674716
$ws.Cells.Item(8,19).Value = 0.08813845619140012
$ws.Cells.Item(8,20).Value = 0.08813845619140014

# Row 9: sCs -> Efnb1/Ephb1 -> sCs
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Efnb1"
$ws.Cells.Item(9,3).Value = "Ephb1"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.273219
$ws.Cells.Item(9,8).Value = 6.819656999999999
$ws.Cells.Item(9,9).Value = 0.1152066857836043
$ws.Cells.Item(9,10).Value = 0.1152066857836043
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.461258
$ws.Cells.Item(9,14).Value = 1.383774
$ws.Cells.Item(9,15).Value = 0.2349536349222574
$ws.Cells.Item(9,16).Value = 0.2349536349222574
$ws.Cells.Item(9,17).Value = 1.048540449502
$ws.Cells.Item(9,18).Value = 9.436864045518
$ws.Cells.Item(9,19).Value = 0.02706822959220419
$ws.Cells.Item(9,20).Value = 0.02706822959220419

Write-Output "done"